# Update accident record in row 4 of the accident report worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new values look numeric to Excel's type-inference but must
#     remain stored as text (matching the original inlineStr cells). Force a
#     text number-format before assigning so Excel does not coerce them into
#     real numbers.
$textCoercedCells = @("G4","J4","Q4","R4","AB4","AC4","AD4","AG4")
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Plain text fields (no numeric-looking coercion risk) ---
$ws.Range("A4").Value = "Head-on collision"
$ws.Range("C4").Value = "2.08 seconds"
$ws.Range("H4").Value = "Car 2"
$ws.Range("S4").Value = "severe"
$ws.Range("X4").Value = "85.3% of historical accidents in Main Highway share severe severity."
$ws.Range("Z4").Value = "severe"
$ws.Range("AA4").Value = "severe"
$ws.Range("AE4").Value = "0.00 (Baseline 10.00, Intervention 10.00)"

# --- Text fields that look numeric (text number-format applied above) ---
$ws.Range("G4").Value = "83.6"
$ws.Range("J4").Value = "-80.8"
$ws.Range("Q4").Value = "2.08"
$ws.Range("R4").Value = "3213773.80"
$ws.Range("AB4").Value = "3213773.80"
$ws.Range("AC4").Value = "2607714.87"
$ws.Range("AD4").Value = "18.86%"
$ws.Range("AG4").Value = "10.00"

# --- Multi-line recommendations text ---
$ws.Range("V4").Value = "- Solar-powered lighting solutions`n- Public awareness campaigns on drinking and driving`n- Install anti-skid road surfaces`n- One-way traffic in narrow roads`n- Install breathalyzer devices"

# --- True numeric cells ---
$ws.Range("F4").Value = 1316
$ws.Range("I4").Value = 1451
$ws.Range("P4").Value = 180

# --- Cells that become blank (previously held Pedestrian info, now unused
#     because this accident involves two vehicles instead of a pedestrian) ---
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
